# TestScenario_eCommerce.xlsx edit
# Commit: "Changed From To postcode"
#
# - Fix the duplicated "TS_006" scenario ID used for both the Checkout
#   section header (B35) and the General GUI verification section
#   (B51:B59); renumber B51:B59 up by one (TS_006->TS_007 ... TS_014->TS_015),
#   introducing a brand-new TS_015 id for the last row.
# - Correct the "No. of test cases" count for the Checkout scenario
#   (G35) from 11 to 15 (there are actually 15 description rows in
#   B35:B49).
# - Update the selected range / active cell shown in the sheet view.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Fix the "No. of test cases" count for the Checkout scenario ---
$ws.Range("G35").Value = 15

# --- Renumber the duplicate TS_006..TS_014 id block to TS_007..TS_015 ---
$ws.Range("B51").Value = "TS_007"
$ws.Range("B52").Value = "TS_008"
$ws.Range("B53").Value = "TS_009"
$ws.Range("B54").Value = "TS_010"
$ws.Range("B55").Value = "TS_011"
$ws.Range("B56").Value = "TS_012"
$ws.Range("B57").Value = "TS_013"
$ws.Range("B58").Value = "TS_014"
$ws.Range("B59").Value = "TS_015"

# --- Update selection / scroll position to reflect where the edit was made ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 31
$excel.ActiveWindow.ScrollColumn = 1
$excel.Goto($ws.Range("C35:C49"))
